$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B6").Value = 0.8090597117364448
$ws.Range("C6").Value = 0.04244249481811561
$ws.Range("E6").Value = 0.7526881720430108
$ws.Range("F6").Value = 0.8172043010752689
$ws.Range("G6").Value = 0.8064516129032258
$ws.Range("H6").Value = 0.8817204301075269
$ws.Range("I6").Value = 0.8091054678563259
$ws.Range("J6").Value = 0.03221425817189315
$ws.Range("M6").Value = 0.8172043010752689
$ws.Range("N6").Value = 0.8172043010752689
$ws.Range("O6").Value = 0.8602150537634409
$ws.Range("P6").Value = 0.7510638297872342
$ws.Range("Q6").Value = 0.0369197168818485
$ws.Range("S6").Value = 0.7204301075268817
$ws.Range("T6").Value = 0.7849462365591398
$ws.Range("U6").Value = 0.6989247311827957
$ws.Range("V6").Value = 0.7956989247311828
$ws.Range("W6").Value = 0.7983299016243424
$ws.Range("X6").Value = 0.02142832910285819
$ws.Range("Y6").Value = 0.776595744680851
$ws.Range("Z6").Value = 0.7849462365591398
$ws.Range("AA6").Value = 0.7956989247311828
$ws.Range("AB6").Value = 0.7956989247311828
$ws.Range("AC6").Value = 0.8387096774193549
$ws.Range("AD6").Value = 0.7983070235644018
$ws.Range("AE6").Value = 0.02058547441921889
$ws.Range("AF6").Value = 0.7872340425531915
$ws.Range("AG6").Value = 0.7849462365591398
$ws.Range("AI6").Value = 0.7741935483870968
$ws.Range("AJ6").Value = 0.8279569892473119
$ws.Range("AK6").Value = 0.8176618622740792
$ws.Range("AL6").Value = 0.04730195882829611
$ws.Range("AM6").Value = 0.7872340425531915
$ws.Range("AN6").Value = 0.7526881720430108
$ws.Range("AO6").Value = 0.8602150537634409
$ws.Range("AP6").Value = 0.8064516129032258
$ws.Range("B7").Value = 0.8541066117593228
$ws.Range("C7").Value = 0.0421345050422935
$ws.Range("D7").Value = 0.8404255319148937
$ws.Range("H7").Value = 0.9139784946236559
$ws.Range("I7").Value = 0.8734156943491191
$ws.Range("J7").Value = 0.04471964960591587
$ws.Range("K7").Value = 0.8617021276595744
$ws.Range("M7").Value = 0.9139784946236559
$ws.Range("N7").Value = 0.8387096774193549
$ws.Range("P7").Value = 0.8326698695950583
$ws.Range("Q7").Value = 0.04089138188946771
$ws.Range("S7").Value = 0.8064516129032258
$ws.Range("T7").Value = 0.8172043010752689
$ws.Range("W7").Value = 0.8798215511324641
$ws.Range("X7").Value = 0.05330075394402922
$ws.Range("Y7").Value = 0.8829787234042553
$ws.Range("Z7").Value = 0.7849462365591398
$ws.Range("AC7").Value = 0.9247311827956989
$ws.Range("AE7").Value = 0.03225874598305675
$ws.Range("AG7").Value = 0.8172043010752689
$ws.Range("AJ7").Value = 0.8817204301075269
$ws.Range("AK7").Value = 0.8647449096316633
$ws.Range("AL7").Value = 0.04850934267717616
$ws.Range("AM7").Value = 0.8936170212765957
$ws.Range("AP7").Value = 0.8602150537634409
